$wb = $excel.ActiveWorkbook

# --- Sheet "股票" (stock): add a "category" column after property_category,
#     and append "source_file" / "index" columns at the end. -----------------
$ws = $wb.Worksheets.Item("股票")

# Insert a new column I ("category"); everything from I onward shifts right.
$ws.Columns("I").Insert()

# Copy the (bold/bordered) header style from the existing "legislator_id"
# header cell (K1) onto the two brand-new trailing header cells so they
# match the rest of the header row instead of falling back to default
# formatting.
$ws.Range("K1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmpe6421"
    $idxVal = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 14).Value = $idxVal
}

# --- Sheet "保險" (insurance): no data changes -- only the shared-string
#     table shifts because of the new strings added above, which Excel
#     handles automatically on save. -----------------------------------------
